$d = $word.ActiveDocument

# --- 1. Remove the three leading navigation paragraphs
#        ("Home", "<- Back to Home", "Download Word Document")
#        that preceded the bookmarked Heading1 title paragraph.
$firstPara = $d.Paragraphs.Item(1)
$thirdPara = $d.Paragraphs.Item(3)
$navText = $d.Range($firstPara.Range.Start, $thirdPara.Range.End).Text
if ($navText -notmatch "Home" -or $navText -notmatch "Download Word Document") {
    throw "Unexpected leading content, aborting: [$navText]"
}
$navRange = $d.Range($firstPara.Range.Start, $thirdPara.Range.End)
$navRange.Delete()

# --- 2. Give every table an explicit 100% preferred width
#        (w:tblW w:type="pct" w:w="5000") instead of the
#        "auto" width Pandoc originally emitted.
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $t.PreferredWidthType = 2   # wdPreferredWidthPercent
    $t.PreferredWidth = 250     # 250 * 20 = 5000 fiftieths-of-a-percent = 100%
}

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
Write-Output "Tables updated: $($d.Tables.Count)"
